$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 538, pushing existing rows 538-641 down to 540-643
$ws.Rows("538:539").Insert()

# Populate new row 538 with the latest weekly entry
$ws.Range("A538").Value = 5
$ws.Range("B538").Value = "Macroferia Regional de Talca"
$ws.Range("C538").Value = "Maule"
$ws.Range("D538").Value = 45209
$ws.Range("E538").Value = 7
$ws.Range("F538").Value = 100112032
$ws.Range("G538").Value = "Zapallo italiano"
$ws.Range("H538").Value = "Sin especificar"
$ws.Range("I538").Value = "Primera"
$ws.Range("J538").Value = 200
$ws.Range("K538").Value = 15000
$ws.Range("L538").Value = 15000
$ws.Range("M538").Value = 15000
$ws.Range("N538").Value = "`$/caja 50 unidades"
$ws.Range("O538").Value = "Región de Arica y Parinacota"
$ws.Range("P538").Value = 300
$ws.Range("Q538").Value = 50
$ws.Range("R538").Value = "Hortaliza"

# Populate new row 539 with the second latest weekly entry
$ws.Range("A539").Value = 5
$ws.Range("B539").Value = "Macroferia Regional de Talca"
$ws.Range("C539").Value = "Maule"
$ws.Range("D539").Value = 45209
$ws.Range("E539").Value = 7
$ws.Range("F539").Value = 100112032
$ws.Range("G539").Value = "Zapallo italiano"
$ws.Range("H539").Value = "Sin especificar"
$ws.Range("I539").Value = "Primera"
$ws.Range("J539").Value = 200
$ws.Range("K539").Value = 17000
$ws.Range("L539").Value = 17000
$ws.Range("M539").Value = 17000
$ws.Range("N539").Value = "`$/caja 60 unidades"
$ws.Range("O539").Value = "Región de O'Higgins"
$ws.Range("P539").Value = 283
$ws.Range("Q539").Value = 60
$ws.Range("R539").Value = "Hortaliza"
